$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.030.03"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.520.29"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.66"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.01"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.518.97"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("E11").Value = "  +3.39%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.117.27"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.63"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000182"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.516.22"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.022.17"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.05"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.33"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.67"
$ws.Range("E21").Value = "  -2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.65"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.76"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.661.94"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.62"
$ws.Range("E28").Value = "  +9.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.60"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.32"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.528.10"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.11"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.146"
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.21"
$ws.Range("E37").Value = "  +4.74%  "
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.45"
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.93"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0807"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.820"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.26"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.25"
$ws.Range("E44").Value = "  +4.17%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.98"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.471.78"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.87"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.895"
$ws.Range("E51").Value = "  +3.38%  "
